$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1046.57
$ws.Range("C3").Value = 1073.18
$ws.Range("C4").Value = 1022.86
$ws.Range("C5").Value = 1061.48
$ws.Range("C6").Value = 1061.48
